# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The account-statement table (rows 16-29) is rebuilt: the two workers
# (DALGIS MARIA LOBO LARA / CC 1045731943 and ENDER MARTINEZ DIAZ / CC
# 8373933) are now interleaved row-by-row for each overdue period
# (2104, 2105, 2108, 2109, 2110, 2111, 2112), instead of being grouped
# worker-by-worker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tipoDoc = "CC"

$dalgisDoc    = "1045731943"
$dalgisNombre = "DALGIS MARIA LOBO LARA"

$enderDoc     = "8373933"
$enderNombre  = "ENDER MARTINEZ DIAZ"

# Each entry: row, periodo, dalgisValorMora, dalgisSalario, enderValorMora, enderSalario
$periods = @(
    @{ Row = 16; Periodo = "2104"; DVM = 35112; DSal = 877803; EVM = 60000; ESal = 1500000 },
    @{ Row = 18; Periodo = "2105"; DVM = 35112; DSal = 877803; EVM = 60000; ESal = 1500000 },
    @{ Row = 20; Periodo = "2108"; DVM = 35112; DSal = 877803; EVM = 60000; ESal = 1500000 },
    @{ Row = 22; Periodo = "2109"; DVM = 35112; DSal = 877803; EVM = 60000; ESal = 1500000 },
    @{ Row = 24; Periodo = "2110"; DVM = 35112; DSal = 877803; EVM = 60000; ESal = 1500000 },
    @{ Row = 26; Periodo = "2111"; DVM = 35112; DSal = 877803; EVM = 60000; ESal = 1500000 },
    @{ Row = 28; Periodo = "2112"; DVM = 18726;  DSal = 877803; EVM = 32000; ESal = 1500000 }
)

foreach ($p in $periods) {
    $rDalgis = $p.Row
    $rEnder  = $p.Row + 1

    $ws.Range("B$rDalgis").Value = $tipoDoc
    $ws.Range("C$rDalgis").Value = $dalgisDoc
    $ws.Range("D$rDalgis").Value = $dalgisNombre
    $ws.Range("E$rDalgis").Value = $p.Periodo
    $ws.Range("F$rDalgis").Value = $p.DVM
    $ws.Range("G$rDalgis").Value = $p.DSal

    $ws.Range("B$rEnder").Value = $tipoDoc
    $ws.Range("C$rEnder").Value = $enderDoc
    $ws.Range("D$rEnder").Value = $enderNombre
    $ws.Range("E$rEnder").Value = $p.Periodo
    $ws.Range("F$rEnder").Value = $p.EVM
    $ws.Range("G$rEnder").Value = $p.ESal
}
